# "changing to use RAST Multiple" - add References / Evidence Types columns
# and switch two genome-name cells from numeric IDs to string IDs (RAST
# multi-genome references), per GramDataEdit5.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (columns C and D)
$ws.Range("C1").Value = "References"
$ws.Range("D1").Value = "Evidence Types"

# New sparse data cells
$ws.Range("C3").Value = "a;b;c"
$ws.Range("D8").Value = "d;e;f"

# A10 must be written before A7 so the shared-string table keeps the same
# ordering as the target workbook (227asdf882.1 before 10asdf6370.11).
$ws.Range("A10").Value = "227asdf882.1"
$ws.Range("A7").Value = "10asdf6370.11"

# Give column D an explicit width (stored width 13.5, matching the target
# <col min="4" max="4" width="13.5" .../> entry).
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666

$ws.Range("B15").Select()
